$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.124.78'
$ws.Range("E2").Value = '  -2.32%  '
$ws.Range("D3").Value = '1.574.68'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.01'
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.495'
$ws.Range("E6").Value = '  -3.65%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0608'
$ws.Range("E8").Value = '  -1.72%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.244'
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.53'
$ws.Range("E10").Value = '  -0.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").Value = '1.797.84'
$ws.Range("E12").Value = '  -1.56%  '
$ws.Range("D13").Value = '1.591.92'
$ws.Range("E13").Value = '  -0.56%  '
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("E15").Value = '  -2.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.29'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("D17").Value = '26.129.73'
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("E19").Value = '  +1.33%  '
$ws.Range("E20").Value = '  -0.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '206.77'
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.25'
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("E23").Value = '  -1.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.85'
$ws.Range("E24").Value = '  -1.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.06'
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.96'
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.111'
$ws.Range("E28").Value = '  -2.13%  '
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0506'
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("E31").Value = '  -1.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.19'
$ws.Range("E32").Value = '  -2.20%  '
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").Value = '1.278.75'
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.45'
$ws.Range("E35").Value = '  -0.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.610'
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.47'
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("E38").Value = '  -2.93%  '
$ws.Range("E39").Value = '  -6.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.814'
$ws.Range("E40").Value = '  -2.33%  '
$ws.Range("E41").Value = '  +2.80%  '
$ws.Range("E42").Value = '  -2.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.762'
$ws.Range("E43").Value = '  -2.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.39'
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("D45").Value = '1.709.98'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.09'
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("D47").Value = '0.0₆0105'
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.52'
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("E50").Value = '  -1.95%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.75'
$ws.Range("E51").Value = '  +10.83%  '
